# Add files via upload — append the new WGAN-P clipping-experiment block
# (rows 18-29) below the existing data on 工作表1, and move the selection
# to the newly-added cell E29.
#
# Values are written in the same order the strings first appear in the
# final sharedStrings table (hyper-parameter block, then the per-layer
# block, then the free-form notes/header, then the row added last) so the
# new shared-string entries line up with the upstream file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 22-24: hyper-parameter block (epoch / critic / clip) ---
$ws.Range("A22").Value = "epoch"
$ws.Range("B22").Value = 10
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 10

$ws.Range("A23").Value = "critic"
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 4

$ws.Range("A24").Value = "clip"
$ws.Range("B24").Value = 0.1
$ws.Range("C24").Value = 0.01
$ws.Range("D24").Value = 0.001

# --- Rows 26-28: per-layer results block ---
$ws.Range("A26").Value = "layer1"
$ws.Range("B26").Value = 4.4956005885581796
$ws.Range("C26").Value = 3.46194675873631
$ws.Range("D26").Value = -9.1322402802024492
$ws.Range("E26").Value = 2.3149222731605801

$ws.Range("A27").Value = "layer2"
$ws.Range("B27").Value = 3.13366123357082
$ws.Range("C27").Value = 2.2632069119990099
$ws.Range("D27").Value = -10.148668956656399
$ws.Range("E27").Value = 1.03625779573425

$ws.Range("A28").Value = "layer3"
$ws.Range("B28").Value = 2.6951344454883599
$ws.Range("C28").Value = 0.99356715988108402
$ws.Range("D28").Value = -10.872890734491
$ws.Range("E28").Value = 0.55493018861950005

# --- Row 25: section label ---
$ws.Range("A25").Value = "base e"

# --- Rows 18-19: free-form Traditional-Chinese observation notes ---
$ws.Range("A18").Value = "應該是c=0.1是越來越小"
$ws.Range("A19").Value = "C=0.01和0.001是越來越大"

# --- Row 21: header row for the new WGAN / WGAN-P comparison table ---
$ws.Range("B21").Value = "WGAN"
$ws.Range("C21").Value = "WGAN"
$ws.Range("D21").Value = "WGAN"
$ws.Range("E21").Value = "WGAN-P"

# --- Row 29: last row added ---
$ws.Range("A29").Value = "layer4"
$ws.Range("B29").Value = 1.010129958142
$ws.Range("C29").Value = -1.08639304558274
$ws.Range("D29").Value = -12.137478965246601
$ws.Range("E29").Value = 0.789269517643991

# --- View state: move selection to the last cell touched ---
$ws.Range("E29").Select()
